$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.043799996376038
$ws.Range("B1").Value = 2.250301599502563
$ws.Range("C1").Value = 4.266026020050049
$ws.Range("D1").Value = 0.8609573245048523
$ws.Range("E1").Value = 1.153002262115479
